# Ajustes de falhas no formulario
# Adds a new data row (row 3) to the "entrada" sheet:
#   1003 | ALCOOL 5L | COZINHA | 2025-03-06 | 45 | COMPRA | 45 | 2025.00 | NF° 45,For:QSUPER

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells -------------------------------------------------
# These strings are not number-like, so a normal .Value assignment is
# already stored as text (matches the style-less text cells already used
# by row 2 for its own text columns).
$ws.Range("B3").Value = "ÁLCOOL 5L"
$ws.Range("C3").Value = "COZINHA"
$ws.Range("F3").Value = "COMPRA"
$ws.Range("I3").Value = "NF° 45,For:QSUPER "

# --- Numeric-looking text cells ---------------------------------------
# "1003", "45", "45" and "2025.00" look like numbers, and a direct
# .Value assignment would store them as real numbers (wrong type, and it
# would also pick up a brand-new number-format style). Row 2 stores the
# analogous columns as plain text with the workbook's default style, so
# we build each value as a text formula in an unused scratch cell (a
# formula result that returns a string is always typed as text,
# regardless of number format) and paste-special *values only* into the
# destination cell. That keeps the destination's style untouched while
# still landing a literal text value.
$ws.Range("K1").Formula = "=""1003"""
$ws.Range("K1").Copy()
$ws.Range("A3").PasteSpecial(-4163)

$ws.Range("K1").Formula = "=""45"""
$ws.Range("K1").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("K1").Formula = "=""45"""
$ws.Range("K1").Copy()
$ws.Range("G3").PasteSpecial(-4163)

$ws.Range("K1").Formula = "=""2025.00"""
$ws.Range("K1").Copy()
$ws.Range("H3").PasteSpecial(-4163)

# Remove the scratch cell completely (value + format) so it leaves no trace.
$ws.Range("K1").Clear()

# --- Date cell ----------------------------------------------------------
# D3 must carry the same date style already used by D2 (numFmtId yyyy-mm-dd).
# Copy D2's formatting only (not its value) so the existing style index is
# reused instead of a new one being created, then write the new date value
# (2025-03-06 -> serial 45722).
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = 45722

$excel.CutCopyMode = 0

Write-Output "Row 3 written"
